$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crawl_time (column F) for all data rows 2-200 to the new scrape timestamp
for ($r = 2; $r -le 200; $r++) {
    $ws.Cells.Item($r, 6).Value = "2025-09-03 09:53:06"
}

# Update date_posted (column B) for rows whose "N hours/days ago" bucket advanced
$bRows = @(2,46,47,48,49,50,51,52,53,76,77,78,79,80,81,84,85,86,87,88,104,105,106,107)
foreach ($r in $bRows) {
    $cell = $ws.Cells.Item($r, 2)
    $old = $cell.Value()
    if ($old -match "^(\d+) (hour|day)s? ago$") {
        $num = [int]$matches[1]
        $unit = $matches[2]
        $newNum = $num + 1
        if ($newNum -eq 1) {
            $suffix = ""
        } else {
            $suffix = "s"
        }
        $cell.Value = "$newNum $unit$suffix ago"
    }
}
